$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Priority" column previously held text values imported from the old
# HANA-based extract (PRIORITY1/2/3). Now that the data comes from SQL
# Server, the priority is a plain numeric column, so replace those text
# cells with numeric 0 placeholders.
$ws.Range("B2:B4").Value = 0

# Leave the selection where the editor last left it.
$ws.Range("D5").Select()
